$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 7 de Abril de 2020 a las 11:22"

$ws.Cells.Item(13, 1).Value = "Belgica"
$ws.Cells.Item(13, 2).Value = 22194
$ws.Cells.Item(13, 3).Value = 1380
$ws.Cells.Item(13, 4).Value = 4157
$ws.Cells.Item(13, 5).Value = 16002
$ws.Cells.Item(13, 6).Value = 1260
$ws.Cells.Item(13, 7).Value = 403
$ws.Cells.Item(13, 8).Value = 2035

$ws.Cells.Item(14, 1).Value = "Suiza"
$ws.Cells.Item(14, 2).Value = 21793
$ws.Cells.Item(14, 3).Value = 136
$ws.Cells.Item(14, 4).Value = 8056
$ws.Cells.Item(14, 5).Value = 12955
$ws.Cells.Item(14, 6).Value = 391
$ws.Cells.Item(14, 7).Value = 17
$ws.Cells.Item(14, 8).Value = 782

$ws.Cells.Item(17, 2).Value = 12399
$ws.Cells.Item(17, 3).Value = 102
$ws.Cells.Item(17, 4).Value = 4046
$ws.Cells.Item(17, 5).Value = 8110
$ws.Cells.Item(17, 7).Value = 23
$ws.Cells.Item(17, 8).Value = 243

$ws.Cells.Item(27, 1).Value = "Dinamarca"
$ws.Cells.Item(27, 2).Value = 4978
$ws.Cells.Item(27, 3).Value = 297
$ws.Cells.Item(27, 4).Value = 1378
$ws.Cells.Item(27, 5).Value = 3413
$ws.Cells.Item(27, 6).Value = 144
$ws.Cells.Item(27, 8).Value = 187

$ws.Cells.Item(28, 1).Value = "India"
$ws.Cells.Item(28, 2).Value = 4858
$ws.Cells.Item(28, 3).Value = 80
$ws.Cells.Item(28, 4).Value = 382
$ws.Cells.Item(28, 5).Value = 4340
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 136

$ws.Cells.Item(29, 1).Value = "Chequia"
$ws.Cells.Item(29, 2).Value = 4828
$ws.Cells.Item(29, 3).Value = 6
$ws.Cells.Item(29, 4).Value = 127
$ws.Cells.Item(29, 5).Value = 4621
$ws.Cells.Item(29, 6).Value = 86
$ws.Cells.Item(29, 7).Value = 2
$ws.Cells.Item(29, 8).Value = 80

$ws.Cells.Item(30, 1).Value = "Chile"
$ws.Cells.Item(30, 2).Value = 4815
$ws.Cells.Item(30, 4).Value = 728
$ws.Cells.Item(30, 5).Value = 4050
$ws.Cells.Item(30, 6).Value = 327
$ws.Cells.Item(30, 8).Value = 37

$ws.Cells.Item(32, 5).Value = 3469
$ws.Cells.Item(32, 7).Value = 6
$ws.Cells.Item(32, 8).Value = 182

$ws.Cells.Item(40, 1).Value = "Indonesia"
$ws.Cells.Item(40, 2).Value = 2738
$ws.Cells.Item(40, 3).Value = 247
$ws.Cells.Item(40, 4).Value = 204
$ws.Cells.Item(40, 5).Value = 2313
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 12
$ws.Cells.Item(40, 8).Value = 221

$ws.Cells.Item(41, 1).Value = "Peru"
$ws.Cells.Item(41, 2).Value = 2561
$ws.Cells.Item(41, 4).Value = 997
$ws.Cells.Item(41, 5).Value = 1472
$ws.Cells.Item(41, 6).Value = 89
$ws.Cells.Item(41, 8).Value = 92

$ws.Cells.Item(71, 1).Value = "Kuwait"
$ws.Cells.Item(71, 2).Value = 743
$ws.Cells.Item(71, 3).Value = 78
$ws.Cells.Item(71, 4).Value = 105
$ws.Cells.Item(71, 5).Value = 637
$ws.Cells.Item(71, 6).Value = 23
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 1

$ws.Cells.Item(72, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(72, 2).Value = 740
$ws.Cells.Item(72, 3).Value = 66
$ws.Cells.Item(72, 4).Value = 68
$ws.Cells.Item(72, 5).Value = 642
$ws.Cells.Item(72, 6).Value = 4
$ws.Cells.Item(72, 7).Value = 1
$ws.Cells.Item(72, 8).Value = 30

$ws.Cells.Item(73, 1).Value = "Crucero"
$ws.Cells.Item(73, 2).Value = 712
$ws.Cells.Item(73, 4).Value = 619
$ws.Cells.Item(73, 5).Value = 82
$ws.Cells.Item(73, 6).Value = 0

$ws.Cells.Item(74, 1).Value = "Bielorrusia"
$ws.Cells.Item(74, 2).Value = 700
$ws.Cells.Item(74, 3).Value = 0
$ws.Cells.Item(74, 4).Value = 53
$ws.Cells.Item(74, 5).Value = 634
$ws.Cells.Item(74, 6).Value = 11
$ws.Cells.Item(74, 8).Value = 13

$ws.Cells.Item(75, 1).Value = "Kazajistan"
$ws.Cells.Item(75, 2).Value = 685
$ws.Cells.Item(75, 3).Value = 23
$ws.Cells.Item(75, 4).Value = 47
$ws.Cells.Item(75, 5).Value = 632
$ws.Cells.Item(75, 6).Value = 16
$ws.Cells.Item(75, 8).Value = 6

$ws.Cells.Item(81, 5).Value = 530
$ws.Cells.Item(81, 7).Value = 1
$ws.Cells.Item(81, 8).Value = 2

$ws.Cells.Item(90, 2).Value = 383
$ws.Cells.Item(90, 3).Value = 6
$ws.Cells.Item(90, 4).Value = 131
$ws.Cells.Item(90, 5).Value = 231

$ws.Cells.Item(116, 1).Value = "Banglades"
$ws.Cells.Item(116, 2).Value = 164
$ws.Cells.Item(116, 3).Value = 41
$ws.Cells.Item(116, 4).Value = 1
$ws.Cells.Item(116, 8).Value = 17

$ws.Cells.Item(117, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(117, 2).Value = 161
$ws.Cells.Item(117, 4).Value = 5
$ws.Cells.Item(117, 5).Value = 138
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 8).Value = 18

$ws.Cells.Item(118, 1).Value = "Kenia"
$ws.Cells.Item(118, 2).Value = 158
$ws.Cells.Item(118, 4).Value = 4
$ws.Cells.Item(118, 5).Value = 148
$ws.Cells.Item(118, 6).Value = 2
$ws.Cells.Item(118, 8).Value = 6

$ws.Cells.Item(119, 1).Value = "Martinica"
$ws.Cells.Item(119, 2).Value = 151
$ws.Cells.Item(119, 4).Value = 50
$ws.Cells.Item(119, 5).Value = 97
$ws.Cells.Item(119, 6).Value = 20
$ws.Cells.Item(119, 8).Value = 4

$ws.Cells.Item(120, 1).Value = "Guadalupe"
$ws.Cells.Item(120, 4).Value = 31
$ws.Cells.Item(120, 5).Value = 101
$ws.Cells.Item(120, 6).Value = 14
$ws.Cells.Item(120, 8).Value = 7

$ws.Cells.Item(121, 1).Value = "Isla de Man"
$ws.Cells.Item(121, 2).Value = 139
$ws.Cells.Item(121, 4).Value = 55
$ws.Cells.Item(121, 5).Value = 83
$ws.Cells.Item(121, 6).Value = 0

$ws.Cells.Item(122, 1).Value = "Brunei"
$ws.Cells.Item(122, 2).Value = 135
$ws.Cells.Item(122, 4).Value = 82
$ws.Cells.Item(122, 5).Value = 52
$ws.Cells.Item(122, 6).Value = 3
$ws.Cells.Item(122, 8).Value = 1

$ws.Cells.Item(123, 1).Value = "Guinea"
$ws.Cells.Item(123, 2).Value = 5
$ws.Cells.Item(123, 6).Value = 0
$ws.Cells.Item(123, 8).Value = 0

$ws.Cells.Item(151, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(151, 3).Value = 15
$ws.Cells.Item(151, 4).Value = 0
$ws.Cells.Item(151, 5).Value = 33
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 8).Value = 0

$ws.Cells.Item(152, 1).Value = "Bahamas"
$ws.Cells.Item(152, 2).Value = 33
$ws.Cells.Item(152, 4).Value = 5
$ws.Cells.Item(152, 5).Value = 23
$ws.Cells.Item(152, 6).Value = 1
$ws.Cells.Item(152, 8).Value = 5

$ws.Cells.Item(153, 1).Value = "Guam"
$ws.Cells.Item(153, 4).Value = 0
$ws.Cells.Item(153, 5).Value = 31
$ws.Cells.Item(153, 6).Value = 0
$ws.Cells.Item(153, 8).Value = 1

$ws.Cells.Item(154, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(154, 2).Value = 32
$ws.Cells.Item(154, 4).Value = 7
$ws.Cells.Item(154, 5).Value = 23
$ws.Cells.Item(154, 6).Value = 6
$ws.Cells.Item(154, 8).Value = 2

$ws.Cells.Item(155, 1).Value = "Eritrea"
$ws.Cells.Item(155, 4).Value = 0
$ws.Cells.Item(155, 5).Value = 31
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 0

$ws.Cells.Item(156, 1).Value = "Guyana"
$ws.Cells.Item(156, 2).Value = 31
$ws.Cells.Item(156, 3).Value = 0
$ws.Cells.Item(156, 4).Value = 8
$ws.Cells.Item(156, 5).Value = 18
$ws.Cells.Item(156, 6).Value = 8
$ws.Cells.Item(156, 7).Value = 1
$ws.Cells.Item(156, 8).Value = 5

$ws.Cells.Item(157, 1).Value = "Gabon"
$ws.Cells.Item(157, 2).Value = 30
$ws.Cells.Item(157, 3).Value = 6
$ws.Cells.Item(157, 4).Value = 1
$ws.Cells.Item(157, 5).Value = 28

$ws.Cells.Item(158, 1).Value = "Benin"
$ws.Cells.Item(158, 2).Value = 26
$ws.Cells.Item(158, 4).Value = 5
$ws.Cells.Item(158, 5).Value = 20

$ws.Cells.Item(159, 1).Value = "Haiti"
$ws.Cells.Item(159, 4).Value = 0
$ws.Cells.Item(159, 5).Value = 23

$ws.Cells.Item(160, 1).Value = "Tanzania"
$ws.Cells.Item(160, 2).Value = 24
$ws.Cells.Item(160, 4).Value = 3
$ws.Cells.Item(160, 5).Value = 20

$ws.Cells.Item(161, 1).Value = "Birmania"
$ws.Cells.Item(161, 2).Value = 22
$ws.Cells.Item(161, 4).Value = 0
$ws.Cells.Item(161, 5).Value = 17

$ws.Cells.Item(162, 1).Value = "Libia"
$ws.Cells.Item(162, 4).Value = 1
$ws.Cells.Item(162, 5).Value = 17
$ws.Cells.Item(162, 8).Value = 1

$ws.Cells.Item(163, 1).Value = "Siria"
$ws.Cells.Item(163, 4).Value = 2
$ws.Cells.Item(163, 5).Value = 15
$ws.Cells.Item(163, 8).Value = 2

$ws.Cells.Item(164, 1).Value = "Maldivas"
$ws.Cells.Item(164, 2).Value = 19
$ws.Cells.Item(164, 4).Value = 13
$ws.Cells.Item(164, 5).Value = 6

$ws.Cells.Item(184, 1).Value = "Zimbabue"
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 8).Value = 1

$ws.Cells.Item(185, 1).Value = "Mozambique"
$ws.Cells.Item(185, 4).Value = 1
$ws.Cells.Item(185, 8).Value = 0

$ws.Cells.Item(193, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(193, 4).Value = 1
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 8).Value = 0

$ws.Cells.Item(195, 1).Value = "Belice"
$ws.Cells.Item(195, 4).Value = 0
$ws.Cells.Item(195, 6).Value = 1
$ws.Cells.Item(195, 8).Value = 1

$ws.Cells.Item(198, 1).Value = "Botsuana"
$ws.Cells.Item(198, 4).Value = 0
$ws.Cells.Item(198, 8).Value = 1

$ws.Cells.Item(200, 1).Value = "San Bartolome"
$ws.Cells.Item(200, 4).Value = 1
$ws.Cells.Item(200, 8).Value = 0

$ws.Cells.Item(204, 1).Value = "Sahara Occidental"

$ws.Cells.Item(205, 1).Value = "Santo Tome y Principe"
